$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.04"
$ws.Range("E2").Value = "'-1.31%"

$ws.Range("D3").Value = "'35.70"
$ws.Range("E3").Value = "'-1.33%"

$ws.Range("D4").Value = "'5.040"
$ws.Range("E4").Value = "'-1.08%"

$ws.Range("D5").Value = "'0.07901"
$ws.Range("E5").Value = "'-2.75%"

$ws.Range("D6").Value = "'1.844"
$ws.Range("E6").Value = "'-5.41%"

$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.101"
$ws.Range("E7").Value = "'-2.13%"

$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.785"
$ws.Range("E8").Value = "'0.25%"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9196"
$ws.Range("E9").Value = "'-1.24%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1353"
$ws.Range("E10").Value = "'-4.61%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1892"
$ws.Range("E11").Value = "'-1.73%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09037"
$ws.Range("E12").Value = "'-2.38%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03471"
$ws.Range("E13").Value = "'-1.39%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09819"
$ws.Range("E14").Value = "'-0.38%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001401"
$ws.Range("E15").Value = "'-0.91%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006116"
$ws.Range("E16").Value = "'4.43%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.723"
$ws.Range("E17").Value = "'3.29%"

$ws.Range("E18").Value = "'10.74%"

$ws.Range("E19").Value = "'-0.01%"

$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1309"
$ws.Range("E20").Value = "'-2.91%"

$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "'5.153"
$ws.Range("E21").Value = "'5.58%"

$ws.Range("E22").Value = "'-8.97%"

$ws.Range("D23").Value = "'0.04403"
$ws.Range("E23").Value = "'-2.39%"

$ws.Range("D24").Value = "'0.001236"
$ws.Range("E24").Value = "'1.44%"

$ws.Range("D25").Value = "'0.004609"
$ws.Range("E25").Value = "'-5.50%"

$ws.Range("E26").Value = "'4.86%"

$ws.Range("D27").Value = "'0.0004444"
$ws.Range("E27").Value = "'0.05%"

$ws.Range("D39").Value = "'0.01935"
$ws.Range("E39").Value = "'-3.72%"

$ws.Range("D40").Value = "'0.05159"
$ws.Range("E40").Value = "'4.55%"

$ws.Range("D41").Value = "'0.007617"
$ws.Range("E41").Value = "'-0.48%"

$ws.Range("D42").Value = "'0.01014"
$ws.Range("E42").Value = "'-5.55%"

$ws.Range("D43").Value = "'0.1339"
$ws.Range("E43").Value = "'-3.16%"

$ws.Range("D44").Value = "'0.002161"
$ws.Range("E44").Value = "'2.88%"

$ws.Range("D45").Value = "'0.01017"
$ws.Range("E45").Value = "'1.75%"

$ws.Range("D46").Value = "'0.00006149"
$ws.Range("E46").Value = "'-4.46%"

$ws.Range("E47").Value = "'0.03%"

$ws.Range("D48").Value = "'63.57"
$ws.Range("E48").Value = "'-1.69%"

$ws.Range("D49").Value = "'0.001660"
$ws.Range("E49").Value = "'39.42%"

$ws.Range("E50").Value = "'0.03%"

$ws.Range("E51").Value = "'0.03%"
